Write-Host "Test1"
Get-ChildItem /tmp/work | Out-String | Write-Host
